# Update countries & provincias Spain
#
# 1. "Last updated" timestamp text changes from 19:20 to 19:50.
# 2. Monaco's case counters are updated and it moves up the (descending,
#    sorted-by-"Casos totales") table, swapping places with Aruba and
#    Banglades (rows 130-132).
# 3. A handful of other countries (Pakistan row 35, Sudafrica row 45,
#    Panama row 48) get small numeric corrections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / "last updated" banner -------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Abril de 2020 a las 19:50"

# --- Pakistan (row 35) ----------------------------------------------------
$ws.Cells.Item(35, 2).Value = 2637
$ws.Cells.Item(35, 3).Value = 216
$ws.Cells.Item(35, 5).Value = 2471

# --- Sudafrica (row 45) ----------------------------------------------------
$ws.Cells.Item(45, 5).Value = 1401
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = 9

# --- Panama (row 48) -------------------------------------------------------
$ws.Cells.Item(48, 4).Value = 10
$ws.Cells.Item(48, 5).Value = 1428

# --- Monaco / Aruba / Banglades re-sort (rows 130-132) ---------------------
# Monaco's "Casos totales" rises to 64, overtaking Aruba (62) and
# Banglades (61); the table stays sorted descending on column B, so
# Monaco now occupies row 130, Aruba row 131 and Banglades row 132.

# Row 130 -> Monaco (new counts)
$ws.Cells.Item(130, 1).Value = "Monaco"
$ws.Cells.Item(130, 2).Value = 64
$ws.Cells.Item(130, 3).Value = 4
$ws.Cells.Item(130, 4).Value = 3
$ws.Cells.Item(130, 5).Value = 60
$ws.Cells.Item(130, 6).Value = 2
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 1

# Row 131 -> Aruba (unchanged counts, shifted down one row)
$ws.Cells.Item(131, 1).Value = "Aruba"
$ws.Cells.Item(131, 2).Value = 62
$ws.Cells.Item(131, 3).Value = 2
$ws.Cells.Item(131, 4).Value = 1
$ws.Cells.Item(131, 5).Value = 61
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 0

# Row 132 -> Banglades (unchanged counts, shifted down one row)
$ws.Cells.Item(132, 1).Value = "Banglades"
$ws.Cells.Item(132, 2).Value = 61
$ws.Cells.Item(132, 3).Value = 5
$ws.Cells.Item(132, 4).Value = 26
$ws.Cells.Item(132, 5).Value = 29
$ws.Cells.Item(132, 6).Value = 1
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 6
